$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.426.58'
$ws.Range("E2").Value = '  -3.87%  '
$ws.Range("D3").Value = '3.022.61'
$ws.Range("E3").Value = '  -3.74%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'538.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.44%  '
$ws.Range("D6").Value = "'132.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -10.96%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.016.54'
$ws.Range("E8").Value = '  -3.67%  '
$ws.Range("E9").Value = '  -3.30%  '
$ws.Range("D10").Value = "'6.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -9.33%  '
$ws.Range("E11").Value = '  -3.85%  '
$ws.Range("E12").Value = '  -2.52%  '
$ws.Range("D13").Value = "'34.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.38%  '
$ws.Range("E14").Value = '  -5.28%  '
$ws.Range("D15").Value = '3.511.31'
$ws.Range("E15").Value = '  -3.66%  '
$ws.Range("D16").Value = '62.478.91'
$ws.Range("E16").Value = '  -3.69%  '
$ws.Range("E17").Value = '  -2.69%  '
$ws.Range("D18").Value = '3.020.22'
$ws.Range("E18").Value = '  -3.76%  '
$ws.Range("D19").Value = "'6.49"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.62%  '
$ws.Range("D20").Value = "'474.29"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -10.37%  '
$ws.Range("E21").Value = '  -4.98%  '
$ws.Range("D22").Value = "'0.684"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = "'6.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -6.46%  '
$ws.Range("D24").Value = "'76.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.11%  '
$ws.Range("D25").Value = "'11.96"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.96%  '
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = "'2.65"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -5.00%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = "'8.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -6.61%  '
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").Value = "'1.89"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -11.03%  '
$ws.Range("D31").Value = "'25.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("E32").Value = '  -4.23%  '
$ws.Range("D33").Value = "'59.70"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +13.00%  '
$ws.Range("E34").Value = '  -7.64%  '
$ws.Range("D35").Value = "'505.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -9.54%  '
$ws.Range("D36").Value = "'5.80"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.37%  '
$ws.Range("D37").Value = "'5.04"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.20%  '
$ws.Range("D38").Value = "'0.0391"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -11.15%  '
$ws.Range("D39").Value = '3.021.53'
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("D40").Value = "'0.0773"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.22%  '
$ws.Range("E41").Value = '  -4.10%  '
$ws.Range("D42").Value = "'7.90"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.20%  '
$ws.Range("D43").Value = "'2.55"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -10.16%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = "'0.247"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.51%  '
$ws.Range("D46").Value = "'1.98"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -8.45%  '
$ws.Range("D47").Value = "'118.61"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("D48").Value = "'23.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.01%  '
$ws.Range("E49").Value = '  -3.45%  '
$ws.Range("D50").Value = "'2.34"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +61.39%  '
$ws.Range("D51").Value = '0.0₃0481'
$ws.Range("E51").Value = '  -8.24%  '
